$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.773.20"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.311.15"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'301.25"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'95.43"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").Value = "'34.17"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'18.94"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "'0.0783"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "2.671.99"
$ws.Range("D16").Value = "2.305.97"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "42.719.32"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("D24").Value = "'235.09"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  +14.59%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'32.08"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "'17.58"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "'0.0697"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'21.13"
$ws.Range("E42").Value = "  +16.55%  "
$ws.Range("D43").Value = "1.924.45"
$ws.Range("D44").Value = "'0.0278"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "'10.08"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "2.541.33"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'53.32"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'72.08"
$ws.Range("E51").Value = "  +1.82%  "
